# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.074.45'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '1.789.47'
$ws.Range("E3").Value = '  +1.63%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.56'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4292'
$ws.Range("E7").Value = '  -3.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3621'
$ws.Range("E8").Value = '  -3.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.85'
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07510'
$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.64'
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.144'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.311'
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").Value = '1.805.68'
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.96'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001066'
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06347'
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.15'
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.977'
$ws.Range("E22").Value = '  -3.63%  '
$ws.Range("D23").Value = '28.096.31'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.37'
$ws.Range("E24").Value = '  -2.39%  '
$ws.Range("E25").Value = '  -7.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.89'
$ws.Range("E26").Value = '  +3.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.35'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '2.006.78'
$ws.Range("E28").Value = '  +2.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.184'
$ws.Range("E29").Value = '  -7.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.27'
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.169'
$ws.Range("E31").Value = '  -4.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.771'
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09009'
$ws.Range("E33").Value = '  -2.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.516'
$ws.Range("E34").Value = '  -3.72%  '
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02334'
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.128'
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6475'
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2114'
$ws.Range("E39").Value = '  -3.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06063'
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.188'
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.420'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.903'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.58'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5996'
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.88'
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.986'
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.154'
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06953'
$ws.Range("E51").Value = '  +0.57%  '
